# Split the movie-title run "Near Real-Time Weather Rendering System – Short Film"
# into "Near Real-Time Weather Rendering System –" (existing run, trimmed) followed
# by a brand-new run containing "Newscast".

$d = $word.ActiveDocument

# 1) Trim the old run's text down to the part that stays in place.
$find1 = $d.Content
$ok1 = $find1.Find.Execute(
    "Near Real-Time Weather Rendering System – Short Film",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Near Real-Time Weather Rendering System –", 2)

# 2) Locate the trimmed text again, collapse to its end (right before the
#    closing guillemet "»") and insert the new word there.
$find2 = $d.Content
$ok2 = $find2.Find.Execute("Near Real-Time Weather Rendering System –")
$find2.Collapse(0)
$find2.InsertAfter("Newscast")

# 3) Force Word to materialise the inserted text as its own run (distinct
#    <w:r>) instead of silently merging it into the preceding run: toggling
#    a direct-formatting property off and back on again splits the run
#    while leaving the resulting formatting identical to its neighbour.
$newRun = $d.Range($find2.Start, $find2.End)
$newRun.Font.Bold = $false
$newRun.Font.Bold = $true
